$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.975.29"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "2.285.79"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'112.47"
$ws.Range("E5").Value = "  -3.73%  "
$ws.Range("D6").Value = "'309.67"
$ws.Range("E6").Value = "  +5.91%  "
$ws.Range("D7").Value = "'0.633"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.615"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").Value = "'44.30"
$ws.Range("E10").Value = "  -4.99%  "
$ws.Range("D11").Value = "'0.0927"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").Value = "'55.07"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("E13").Value = "  -3.76%  "
$ws.Range("D14").Value = "'1.09"
$ws.Range("E14").Value = "  +21.31%  "
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").Value = "'15.51"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "2.629.16"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "2.284.13"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("D19").Value = "42.925.59"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("D21").Value = "'7.19"
$ws.Range("E21").Value = "  -6.25%  "
$ws.Range("D22").Value = "'76.52"
$ws.Range("E22").Value = "  +3.63%  "
$ws.Range("D23").Value = "'3.58"
$ws.Range("E23").Value = "  +3.95%  "
$ws.Range("D24").Value = "'2.46"
$ws.Range("E24").Value = "  +3.34%  "
$ws.Range("D25").Value = "'256.86"
$ws.Range("E25").Value = "  +9.97%  "
$ws.Range("D26").Value = "'8.96"
$ws.Range("E26").Value = "  -4.49%  "
$ws.Range("D27").Value = "'11.74"
$ws.Range("E27").Value = "  -4.07%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("D30").Value = "'38.76"
$ws.Range("E30").Value = "  -4.17%  "
$ws.Range("D31").Value = "'22.32"
$ws.Range("E31").Value = "  +4.54%  "
$ws.Range("D32").Value = "'173.83"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("D34").Value = "'0.0899"
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("D36").Value = "'5.07"
$ws.Range("E36").Value = "  +6.52%  "
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").Value = "'4.16"
$ws.Range("E38").Value = "  -9.36%  "
$ws.Range("D39").Value = "'0.0376"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("E40").Value = "  -2.89%  "
$ws.Range("D41").Value = "'2.53"
$ws.Range("D42").Value = "'72.37"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'1.39"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "'12.42"
$ws.Range("E46").Value = "  -8.53%  "
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("D48").Value = "'108.11"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("D49").Value = "'8.89"
$ws.Range("E49").Value = "  +2.69%  "
$ws.Range("D50").Value = "'1.30"
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("D51").Value = "'71.25"
$ws.Range("E51").Value = "  -0.55%  "
